$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "50.720.44"
Set-TextCell 2 5 "  -1.12%  "
Set-TextCell 3 4 "2.918.41"
Set-TextCell 3 5 "  -1.79%  "
Set-TextCell 4 4 "1.00"
Set-TextCell 4 5 "  +0.04%  "
Set-TextCell 5 4 "374.58"
Set-TextCell 5 5 "  -2.03%  "
Set-TextCell 6 4 "99.52"
Set-TextCell 6 5 "  -3.16%  "
Set-TextCell 7 5 "  -1.05%  "
Set-TextCell 8 5 "  +0.00%  "
Set-TextCell 9 4 "0.575"
Set-TextCell 9 5 "  -2.80%  "
Set-TextCell 10 4 "35.56"
Set-TextCell 10 5 "  -2.89%  "
Set-TextCell 11 5 "  -0.76%  "
Set-TextCell 12 4 "0.0844"
Set-TextCell 12 5 "  +0.38%  "
Set-TextCell 13 4 "3.379.78"
Set-TextCell 13 5 "  -1.90%  "
Set-TextCell 14 4 "17.93"
Set-TextCell 14 5 "  -1.11%  "
Set-TextCell 15 4 "7.59"
Set-TextCell 15 5 "  +1.75%  "
Set-TextCell 16 4 "11.95"
Set-TextCell 16 5 "  +62.56%  "
Set-TextCell 17 4 "2.935.39"
Set-TextCell 17 5 "  -0.95%  "
Set-TextCell 18 4 "0.991"
Set-TextCell 18 5 "  +0.07%  "
Set-TextCell 19 4 "50.676.28"
Set-TextCell 19 5 "  -1.07%  "
Set-TextCell 20 5 "  -7.61%  "
Set-TextCell 21 4 "12.22"
Set-TextCell 21 5 "  -3.82%  "
Set-TextCell 22 5 "  -1.61%  "
Set-TextCell 23 4 "69.23"
Set-TextCell 23 5 "  +0.66%  "
Set-TextCell 24 4 "265.29"
Set-TextCell 24 5 "  +1.22%  "
Set-TextCell 25 4 "3.13"
Set-TextCell 25 5 "  +8.40%  "
Set-TextCell 26 4 "7.80"
Set-TextCell 26 5 "  -4.03%  "
Set-TextCell 27 5 "  -0.01%  "
Set-TextCell 28 4 "7.09"
Set-TextCell 28 5 "  -6.18%  "
Set-TextCell 29 4 "25.31"
Set-TextCell 29 5 "  -1.80%  "
Set-TextCell 30 5 "  -3.44%  "
Set-TextCell 31 5 "  -5.08%  "
Set-TextCell 32 5 "  +0.41%  "
Set-TextCell 33 4 "50.26"
Set-TextCell 33 5 "  -1.30%  "
Set-TextCell 34 5 "  -0.17%  "
Set-TextCell 35 4 "33.03"
Set-TextCell 35 5 "  -3.18%  "
Set-TextCell 36 5 "  -3.88%  "
Set-TextCell 37 5 "  -0.08%  "
Set-TextCell 38 5 "  +1.93%  "
Set-TextCell 39 4 "0.115"
Set-TextCell 39 5 "  -0.65%  "
Set-TextCell 40 4 "16.24"
Set-TextCell 40 5 "  -4.80%  "
Set-TextCell 41 5 "  -0.89%  "
Set-TextCell 42 4 "2.40"
Set-TextCell 42 5 "  -6.05%  "
Set-TextCell 43 4 "119.18"
Set-TextCell 43 5 "  -2.63%  "
Set-TextCell 44 4 "20.69"
Set-TextCell 44 5 "  -2.88%  "
Set-TextCell 45 2 "WEMIXToken"
Set-TextCell 45 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 45 4 "2.03"
Set-TextCell 45 5 "  -1.98%  "
Set-TextCell 46 2 "NEARProtocol"
Set-TextCell 46 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 46 4 "3.33"
Set-TextCell 46 5 "  +2.58%  "
Set-TextCell 48 4 "1.984.90"
Set-TextCell 49 4 "0.256"
Set-TextCell 49 5 "  -6.75%  "
Set-TextCell 50 4 "0.0312"
Set-TextCell 50 5 "  -6.85%  "
Set-TextCell 51 4 "5.24"
Set-TextCell 51 5 "  +3.01%  "
